# Refresh the cryptos table: updated Price (D) / Volume(1h) (E) values,
# plus the Stellar/Hedera rows (48-49) swapping rank order.
# Note: D12 and D23 use a leading apostrophe to force literal text, since
# "0.150"/"1.00" would otherwise be auto-converted to numbers (0.15 / 1),
# dropping the significant trailing zero that the source text preserves.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.849.09'
$ws.Range('E2').Value = '  +4.80%  '
$ws.Range('D3').Value = '2.434.15'
$ws.Range('E3').Value = '  +5.31%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '565.23'
$ws.Range('E5').Value = '  +4.30%  '
$ws.Range('D6').Value = '140.48'
$ws.Range('E6').Value = '  +8.10%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.586'
$ws.Range('E8').Value = '  +2.30%  '
$ws.Range('D9').Value = '2.432.41'
$ws.Range('E9').Value = '  +5.33%  '
$ws.Range('E10').Value = '  +3.80%  '
$ws.Range('E11').Value = '  +4.09%  '
$ws.Range('D12').Value = "'0.150"
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('D13').Value = '0.349'
$ws.Range('E13').Value = '  +5.39%  '
$ws.Range('D14').Value = '26.42'
$ws.Range('E14').Value = '  +13.27%  '
$ws.Range('D15').Value = '2.864.98'
$ws.Range('E15').Value = '  +5.23%  '
$ws.Range('D16').Value = '62.737.97'
$ws.Range('E16').Value = '  +4.64%  '
$ws.Range('E17').Value = '  +8.23%  '
$ws.Range('D18').Value = '2.435.19'
$ws.Range('E18').Value = '  +5.84%  '
$ws.Range('D19').Value = '11.26'
$ws.Range('E19').Value = '  +7.55%  '
$ws.Range('D20').Value = '340.53'
$ws.Range('E20').Value = '  +9.09%  '
$ws.Range('D21').Value = '4.22'
$ws.Range('E21').Value = '  +3.56%  '
$ws.Range('D22').Value = '6.81'
$ws.Range('E22').Value = '  +4.05%  '
$ws.Range('D23').Value = "'1.00"
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = '5.64'
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').Value = '65.48'
$ws.Range('E25').Value = '  +3.44%  '
$ws.Range('E26').Value = '  +2.03%  '
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').Value = '1.54'
$ws.Range('E28').Value = '  +14.34%  '
$ws.Range('D29').Value = '8.23'
$ws.Range('E29').Value = '  +6.30%  '
$ws.Range('D30').Value = '1.36'
$ws.Range('E30').Value = '  +15.44%  '
$ws.Range('E31').Value = '  +9.04%  '
$ws.Range('E32').Value = '  +5.99%  '
$ws.Range('D33').Value = '6.55'
$ws.Range('E33').Value = '  +12.22%  '
$ws.Range('D34').Value = '174.08'
$ws.Range('E34').Value = '  +1.47%  '
$ws.Range('E35').Value = '  +8.57%  '
$ws.Range('D36').Value = '0.398'
$ws.Range('E36').Value = '  +5.17%  '
$ws.Range('D37').Value = '378.32'
$ws.Range('E37').Value = '  +19.55%  '
$ws.Range('D38').Value = '18.61'
$ws.Range('E38').Value = '  +5.18%  '
$ws.Range('E39').Value = '  +11.85%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('E42').Value = '  +13.39%  '
$ws.Range('D43').Value = '39.92'
$ws.Range('E43').Value = '  +5.90%  '
$ws.Range('D44').Value = '144.93'
$ws.Range('E44').Value = '  +6.62%  '
$ws.Range('D45').Value = '3.67'
$ws.Range('E45').Value = '  +7.06%  '
$ws.Range('D46').Value = '20.59'
$ws.Range('E46').Value = '  +10.06%  '
$ws.Range('D47').Value = '0.594'
$ws.Range('E47').Value = '  +4.42%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').Value = '0.0521'
$ws.Range('E48').Value = '  +6.47%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '0.0948'
$ws.Range('E49').Value = '  +0.70%  '
$ws.Range('E50').Value = '  +5.13%  '
$ws.Range('D51').Value = '17.89'
$ws.Range('E51').Value = '  +7.07%  '
